$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 175.88889
$ws.Range("I9").Value = 65.14286
$ws.Range("K9").Value = 65.14286
$ws.Range("M9").Value = 103.85714
$ws.Range("H17").Value = 2733414.8
$ws.Range("J17").Value = 2733414.8
$ws.Range("L17").Value = 8200244.399999999
$ws.Range("N17").Value = -8200580.399999999
$ws.Range("H18").Value = 215.5
$ws.Range("I18").Value = 215.5
$ws.Range("K18").Value = 215.5
$ws.Range("M18").Value = 68.5
$ws.Range("H38").Value = 3483.4285
$ws.Range("I38").Value = 390
$ws.Range("K38").Value = 1170
$ws.Range("M38").Value = -798
$ws.Range("H40").Value = 1728.4286
$ws.Range("J40").Value = 1516.5
$ws.Range("L40").Value = 1516.5
$ws.Range("N40").Value = -1866.5
$ws.Range("H48").Value = 2000
$ws.Range("I48").Value = 2000
$ws.Range("K48").Value = 6000
$ws.Range("M48").Value = -5708
$ws.Range("H56").Value = 2000
$ws.Range("I56").Value = 2000
$ws.Range("K56").Value = 6000
$ws.Range("M56").Value = -5466
$ws.Range("H64").Value = 500000000
$ws.Range("I64").Value = 500000000
$ws.Range("K64").Value = 500000000
$ws.Range("M64").Value = -499999752
$ws.Range("H67").Value = 500000000
$ws.Range("I67").Value = 500000000
$ws.Range("K67").Value = 500000000
$ws.Range("M67").Value = -499999142
$ws.Range("H70").Value = 2614.7144
$ws.Range("I70").Value = 2500
$ws.Range("K70").Value = 7500
$ws.Range("M70").Value = -7230
$ws.Range("H73").Value = 2614.7144
$ws.Range("I73").Value = 2500
$ws.Range("K73").Value = 7500
$ws.Range("M73").Value = -6564
$ws.Range("H80").Value = 1006.875
$ws.Range("I80").Value = 982.5
$ws.Range("J80").Value = 1031.25
$ws.Range("K80").Value = 2947.5
$ws.Range("L80").Value = 3093.75
$ws.Range("M80").Value = -1949.5
$ws.Range("N80").Value = -5089.75
$ws.Range("H83").Value = 1006.875
$ws.Range("I83").Value = 982.5
$ws.Range("J83").Value = 1031.25
$ws.Range("K83").Value = 8842.5
$ws.Range("L83").Value = 9281.25
$ws.Range("M83").Value = -3850.5
$ws.Range("N83").Value = -19265.25
$ws.Range("H96").Value = 1143.9615
$ws.Range("I96").Value = 949.2632
$ws.Range("J96").Value = 1672.4286
$ws.Range("K96").Value = 2847.7896
$ws.Range("L96").Value = 5017.2858
$ws.Range("M96").Value = -1474.7896
$ws.Range("N96").Value = -7763.2858
$ws.Range("H98").Value = 2949.2104
$ws.Range("I98").Value = 2620.8823
$ws.Range("K98").Value = 2620.8823
$ws.Range("M98").Value = -1122.8823
$ws.Range("H106").Value = 4487.75
$ws.Range("I106").Value = 4487.75
$ws.Range("K106").Value = 4487.75
$ws.Range("M106").Value = -3856.75
$ws.Range("H122").Value = 2949.2104
$ws.Range("I122").Value = 2620.8823
$ws.Range("K122").Value = 7862.646900000001
$ws.Range("M122").Value = -5412.646900000001
$ws.Range("H125").Value = 3250813.8
$ws.Range("I125").Value = 5055403
$ws.Range("J125").Value = 2552.8
$ws.Range("K125").Value = 45498627
$ws.Range("L125").Value = 22975.2
$ws.Range("M125").Value = -45496167
$ws.Range("N125").Value = -27895.2
$ws.Range("H132").Value = 18828.137
$ws.Range("I132").Value = 11906.353
$ws.Range("K132").Value = 35719.05899999999
$ws.Range("M132").Value = -33189.05899999999
$ws.Range("H137").Value = 18897
$ws.Range("I137").Value = 1807.8572
$ws.Range("K137").Value = 5423.571599999999
$ws.Range("M137").Value = -2873.571599999999
$ws.Range("H138").Value = 2130.1191
$ws.Range("J138").Value = 2351.362
$ws.Range("L138").Value = 7054.086
$ws.Range("N138").Value = -17334.086

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3516.8867
$ws.Range("I32").Value = 2244.4666
$ws.Range("K32").Value = 2244.4666
$ws.Range("M32").Value = -1957.4666
$ws.Range("H33").Value = 3750
$ws.Range("I33").Value = 2500
$ws.Range("K33").Value = 2500
$ws.Range("M33").Value = -2171
$ws.Range("H52").Value = 110000
$ws.Range("J52").Value = 110000
$ws.Range("L52").Value = 110000
$ws.Range("N52").Value = -110636
$ws.Range("H61").Value = 55915.652
$ws.Range("I61").Value = 1339.069
$ws.Range("K61").Value = 1339.069
$ws.Range("M61").Value = -1127.069
$ws.Range("H63").Value = 3050.375
$ws.Range("I63").Value = 3092.1667
$ws.Range("J63").Value = 2925
$ws.Range("K63").Value = 3092.1667
$ws.Range("L63").Value = 2925
$ws.Range("M63").Value = -2406.1667
$ws.Range("N63").Value = -4297
$ws.Range("H66").Value = 3050.375
$ws.Range("I66").Value = 3092.1667
$ws.Range("J66").Value = 2925
$ws.Range("K66").Value = 15460.8335
$ws.Range("L66").Value = 14625
$ws.Range("M66").Value = -12028.8335
$ws.Range("N66").Value = -21489
$ws.Range("H74").Value = 10959.303
$ws.Range("I74").Value = 1574.9459
$ws.Range("J74").Value = 68829.5
$ws.Range("K74").Value = 1574.9459
$ws.Range("L74").Value = 68829.5
$ws.Range("M74").Value = -700.9458999999999
$ws.Range("N74").Value = -70577.5
$ws.Range("H77").Value = 10959.303
$ws.Range("I77").Value = 1574.9459
$ws.Range("J77").Value = 68829.5
$ws.Range("K77").Value = 7874.729499999999
$ws.Range("L77").Value = 344147.5
$ws.Range("M77").Value = -3506.729499999999
$ws.Range("N77").Value = -352883.5
$ws.Range("H110").Value = 5351233
$ws.Range("I110").Value = 5351233
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 5351233
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -5349188
$ws.Range("H136").Value = 55915.652
$ws.Range("I136").Value = 1339.069
$ws.Range("K136").Value = 4017.207
$ws.Range("M136").Value = -1467.207
$ws.Range("N110").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 29412638
$ws.Range("I105").Value = 33334206
$ws.Range("J105").Value = 884.5
$ws.Range("K105").Value = 33334206
$ws.Range("L105").Value = 884.5
$ws.Range("M105").Value = -33332459
$ws.Range("N105").Value = -4378.5
$ws.Range("H134").Value = 32228.75
$ws.Range("I134").Value = 35437.8
$ws.Range("K134").Value = 106313.4
$ws.Range("M134").Value = -103778.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 559.0769
$ws.Range("I7").Value = 145.85715
$ws.Range("J7").Value = 1041.1666
$ws.Range("K7").Value = 145.85715
$ws.Range("L7").Value = 1041.1666
$ws.Range("M7").Value = -32.85714999999999
$ws.Range("N7").Value = -1267.1666
$ws.Range("H22").Value = 1391
$ws.Range("J22").Value = 1841.7142
$ws.Range("L22").Value = 1841.7142
$ws.Range("N22").Value = -2541.7142
$ws.Range("H58").Value = 11789.768
$ws.Range("J58").Value = 25459.285
$ws.Range("L58").Value = 25459.285
$ws.Range("N58").Value = -25865.285
$ws.Range("H99").Value = 8405222
$ws.Range("J99").Value = 6671531.5
$ws.Range("L99").Value = 6671531.5
$ws.Range("N99").Value = -6674527.5
$ws.Range("H125").Value = 85003.664
$ws.Range("J125").Value = 85003.664
$ws.Range("L125").Value = 85003.664
$ws.Range("N125").Value = -89923.664
$ws.Range("H126").Value = 8405222
$ws.Range("J126").Value = 6671531.5
$ws.Range("L126").Value = 20014594.5
$ws.Range("N126").Value = -20019534.5
$ws.Range("H132").Value = 4251
$ws.Range("I132").Value = 3775.3076
$ws.Range("K132").Value = 11325.9228
$ws.Range("M132").Value = -8795.9228
$ws.Range("H133").Value = 106000
$ws.Range("J133").Value = 106000
$ws.Range("L133").Value = 106000
$ws.Range("N133").Value = -111060
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 11789.768
$ws.Range("J136").Value = 25459.285
$ws.Range("L136").Value = 76377.855
$ws.Range("N136").Value = -81477.855
$ws.Range("N135").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2313271
$ws.Range("I4").Value = 3219785.5
$ws.Range("K4").Value = 9659356.5
$ws.Range("M4").Value = -9659244.5
$ws.Range("H5").Value = 11150939
$ws.Range("I5").Value = 1118.8889
$ws.Range("K5").Value = 3356.6667
$ws.Range("M5").Value = -3244.6667
$ws.Range("H68").Value = 4299.25
$ws.Range("J68").Value = 4299.25
$ws.Range("L68").Value = 12897.75
$ws.Range("N68").Value = -14519.75
$ws.Range("H71").Value = 4299.25
$ws.Range("J71").Value = 4299.25
$ws.Range("L71").Value = 38693.25
$ws.Range("N71").Value = -46805.25
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("H114").Value = 1229.909
$ws.Range("I114").Value = 199.66667
$ws.Range("J114").Value = 1616.25
$ws.Range("K114").Value = 599.00001
$ws.Range("L114").Value = 4848.75
$ws.Range("M114").Value = 2654.99999
$ws.Range("N114").Value = -11356.75
$ws.Range("H117").Value = 1410.625
$ws.Range("I117").Value = 247
$ws.Range("J117").Value = 2574.25
$ws.Range("K117").Value = 741
$ws.Range("L117").Value = 7722.75
$ws.Range("M117").Value = 2701
$ws.Range("N117").Value = -14606.75
$ws.Range("H130").Value = 16530.77
$ws.Range("I130").Value = 4966.6665
$ws.Range("J130").Value = 20000
$ws.Range("K130").Value = 14899.9995
$ws.Range("L130").Value = 60000
$ws.Range("M130").Value = -9879.999500000002
$ws.Range("N130").Value = -70040
$ws.Range("H131").Value = 1421.39
$ws.Range("I131").Value = 499
$ws.Range("J131").Value = 1430.707
$ws.Range("K131").Value = 1497
$ws.Range("L131").Value = 4292.121
$ws.Range("M131").Value = 3543
$ws.Range("N131").Value = -14372.121
$ws.Range("H132").Value = 2246.5
$ws.Range("I132").Value = 1955.7778
$ws.Range("K132").Value = 17602.0002
$ws.Range("M132").Value = -15072.0002
$ws.Range("H135").Value = 11150939
$ws.Range("I135").Value = 1118.8889
$ws.Range("K135").Value = 10070.0001
$ws.Range("M135").Value = -7535.000099999999
$ws.Range("M76").ClearContents()
$ws.Range("M79").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 8382
$ws.Range("I21").Value = 3250
$ws.Range("J21").Value = 10948
$ws.Range("K21").Value = 3250
$ws.Range("L21").Value = 10948
$ws.Range("M21").Value = -3077
$ws.Range("N21").Value = -11294
$ws.Range("H30").Value = 8382
$ws.Range("I30").Value = 3250
$ws.Range("J30").Value = 10948
$ws.Range("K30").Value = 3250
$ws.Range("L30").Value = 10948
$ws.Range("M30").Value = -3145
$ws.Range("N30").Value = -11158
$ws.Range("H113").Value = 3378.3
$ws.Range("I113").Value = 2598.8
$ws.Range("K113").Value = 2598.8
$ws.Range("M113").Value = -428.8000000000002
$ws.Range("H122").Value = 6786508
$ws.Range("I122").Value = 6786508
$ws.Range("K122").Value = 20359524
$ws.Range("M122").Value = -20357074
$ws.Range("H126").Value = 7163578
$ws.Range("I126").Value = 4000647
$ws.Range("J126").Value = 11907975
$ws.Range("K126").Value = 12001941
$ws.Range("L126").Value = 35723925
$ws.Range("M126").Value = -11999471
$ws.Range("N126").Value = -35728865
$ws.Range("H132").Value = 2848.0557
$ws.Range("I132").Value = 2578.1875
$ws.Range("K132").Value = 7734.5625
$ws.Range("M132").Value = -5204.5625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1709258.9
$ws.Range("I7").Value = 3183845.2
$ws.Range("K7").Value = 3183845.2
$ws.Range("M7").Value = -3183733.2
$ws.Range("H32").Value = 7589.25
$ws.Range("I32").Value = 1699.75
$ws.Range("J32").Value = 13478.75
$ws.Range("K32").Value = 1699.75
$ws.Range("L32").Value = 13478.75
$ws.Range("M32").Value = -1382.75
$ws.Range("N32").Value = -14112.75
$ws.Range("H40").Value = 1735689.6
$ws.Range("I40").Value = 4562.7617
$ws.Range("J40").Value = 4532125
$ws.Range("K40").Value = 4562.7617
$ws.Range("L40").Value = 4532125
$ws.Range("M40").Value = -4426.7617
$ws.Range("N40").Value = -4532397
$ws.Range("H46").Value = 1855.1111
$ws.Range("I46").Value = 750
$ws.Range("J46").Value = 2170.8572
$ws.Range("K46").Value = 750
$ws.Range("L46").Value = 2170.8572
$ws.Range("M46").Value = -562
$ws.Range("N46").Value = -2546.8572
$ws.Range("H56").Value = 34500
$ws.Range("J56").Value = 50000
$ws.Range("L56").Value = 50000
$ws.Range("N56").Value = -51382
$ws.Range("H59").Value = 50133
$ws.Range("J59").Value = 50133
$ws.Range("L59").Value = 50133
$ws.Range("N59").Value = -51441
$ws.Range("H126").Value = 1709258.9
$ws.Range("I126").Value = 3183845.2
$ws.Range("K126").Value = 9551535.600000001
$ws.Range("M126").Value = -9549065.600000001
$ws.Range("H136").Value = 10990.155
$ws.Range("I136").Value = 8678.9
$ws.Range("J136").Value = 15612.667
$ws.Range("K136").Value = 26036.7
$ws.Range("L136").Value = 46838.001
$ws.Range("M136").Value = -23486.7
$ws.Range("N136").Value = -51938.001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 2966.6667
$ws.Range("I58").Value = 2966.6667
$ws.Range("K58").Value = 2966.6667
$ws.Range("M58").Value = -2658.6667
$ws.Range("H76").Value = 45171.332
$ws.Range("J76").Value = 45171.332
$ws.Range("L76").Value = 45171.332
$ws.Range("N76").Value = -45801.332
$ws.Range("H79").Value = 45171.332
$ws.Range("J79").Value = 45171.332
$ws.Range("L79").Value = 45171.332
$ws.Range("N79").Value = -47355.332
$ws.Range("H81").Value = 2765.4707
$ws.Range("I81").Value = 2914.5625
$ws.Range("J81").Value = 380
$ws.Range("K81").Value = 5829.125
$ws.Range("L81").Value = 760
$ws.Range("M81").Value = -4768.125
$ws.Range("N81").Value = -2882
$ws.Range("H84").Value = 2765.4707
$ws.Range("I84").Value = 2914.5625
$ws.Range("J84").Value = 380
$ws.Range("K84").Value = 29145.625
$ws.Range("L84").Value = 3800
$ws.Range("M84").Value = -23841.625
$ws.Range("N84").Value = -14408
$ws.Range("H97").Value = 26688.857
$ws.Range("J97").Value = 26688.857
$ws.Range("L97").Value = 26688.857
$ws.Range("N97").Value = -28670.857
$ws.Range("H107").Value = 1140.4375
$ws.Range("I107").Value = 1188.2307
$ws.Range("K107").Value = 3564.6921
$ws.Range("M107").Value = -1644.6921
$ws.Range("H126").Value = 6254414
$ws.Range("I126").Value = 5617.25
$ws.Range("J126").Value = 25000804
$ws.Range("K126").Value = 16851.75
$ws.Range("L126").Value = 75002412
$ws.Range("M126").Value = -14381.75
$ws.Range("N126").Value = -75007352
$ws.Range("H136").Value = 15355.454
$ws.Range("I136").Value = 2542.318
$ws.Range("K136").Value = 7626.954000000001
$ws.Range("M136").Value = -5076.954000000001
